# Generate Report for Handback
# Update the "generate date" / handoff / handback timestamp cells that get
# refreshed whenever the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Column G = "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 11:08:57"

# --- zh-cn sheet ------------------------------------------------------
# H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 11:08:52"
$wsZhCn.Range("K2").Value = "2016-08-25 11:09:16"

# --- de-de sheet ------------------------------------------------------
# H2 = Correspond Handoff Datetime (shares the same generate-date value as
# the Overview sheet's G2), K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 11:08:57"
$wsDeDe.Range("K2").Value = "2016-08-25 11:09:24"
